# Add a second row of data below the existing header row.
# Values must stay literal text (e.g. "09/08/2023" / "5000.00"), not be
# coerced into dates/numbers, so force the range to Text format before
# writing, then drop back to the default "Normal" style so no extra
# formatting is left on the new cells (matches the source workbook, whose
# header row is the only one with an explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A2:H2")
$newRow.NumberFormat = "@"

$ws.Range("A2").Value = "09/08/2023"
$ws.Range("B2").Value = "5000.00"
$ws.Range("C2").Value = "5000.00"
$ws.Range("D2").Value = "4000.00"
$ws.Range("E2").Value = "4000.00"
$ws.Range("F2").Value = "20.0"
$ws.Range("G2").Value = "1000.00"
$ws.Range("H2").Value = "80.00"

$newRow.Style = "Normal"
